$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F:V) between rows 99 and 101 ---
# (columns A:E -- index/pais/torneio/temporada/data_partida -- stay put)
$row99 = $ws.Range("F99:V99").Value()
$row101 = $ws.Range("F101:V101").Value()
$ws.Range("F99:V99").Value = $row101
$ws.Range("F101:V101").Value = $row99

# --- Swap the match data (columns F:V) between rows 137 and 138 ---
$row137 = $ws.Range("F137:V137").Value()
$row138 = $ws.Range("F138:V138").Value()
$ws.Range("F137:V137").Value = $row138
$ws.Range("F138:V138").Value = $row137

# --- Swap the match data (columns F:V) between rows 139 and 140 ---
$row139 = $ws.Range("F139:V139").Value()
$row140 = $ws.Range("F140:V140").Value()
$ws.Range("F139:V139").Value = $row140
$ws.Range("F140:V140").Value = $row139

# --- Append the new match as row 141 (Palermo vs Catanzaro) ---
# Copy formatting from the row above first so styles (bold index cell,
# date/time number format, borders, ...) match the rest of the table.
$ws.Range("A140:V140").Copy()
$ws.Range("A141:V141").PasteSpecial(-4122)

$ws.Range("A141").Value = 140
$ws.Range("B141").Value = "italy"
$ws.Range("C141").Value = "serie-b"
$ws.Range("D141").Value = "2023-2024"
$ws.Range("E141").Value = 45261.85416666666
$ws.Range("F141").Value = "Palermo"
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = "Catanzaro"
$ws.Range("I141").Value = 2
$ws.Range("J141").Value = 1.86
$ws.Range("K141").Value = "26/11/2023 16:43"
$ws.Range("L141").Value = 2.1
$ws.Range("M141").Value = "01/12/2023 20:29"
$ws.Range("N141").Value = 3.78
$ws.Range("O141").Value = "26/11/2023 16:43"
$ws.Range("P141").Value = 3.34
$ws.Range("Q141").Value = "01/12/2023 20:29"
$ws.Range("R141").Value = 4.3
$ws.Range("S141").Value = "26/11/2023 16:43"
$ws.Range("T141").Value = 3.95
$ws.Range("U141").Value = "01/12/2023 20:29"
$ws.Range("V141").Value = "https://www.betexplorer.com/football/italy/serie-b/palermo-catanzaro/Ym99AElE/"
